# Add a new "PF/1.0.5" row to the meta-sheet, per commit message
# "Add PF/1.0.5 to meta-sheet".
#
# Existing layout (row 1 = environment headers, row 2 = current prod
# version "PF/1.0.0" repeated across columns). We append row 3 with the
# new version id in column A and "X" markers in columns B:D.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "PF/1.0.5"
$ws.Range("B3").Value = "X"
$ws.Range("C3").Value = "X"
$ws.Range("D3").Value = "X"

# New row keeps the workbook's default (unstyled) formatting instead of
# inheriting the header/body alignment style used by rows 1-2.
$ws.Range("A3:D3").Style = "Normal"
